# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets
# to match the newly scraped data (output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Row -> new value for column F
$updates = @{
    2  = 1075
    3  = 780
    8  = 1912
    9  = 6783
    11 = 389
    12 = 322
    14 = 385
    16 = 7005
    17 = 284
    21 = 225
    28 = 15
    29 = 401
    34 = 67
    35 = 30
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
